$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.770.25"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.635.78"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'215.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "'19.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'0.0786"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "1.634.69"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "1.861.38"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'0.556"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").Value = "'63.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "25.800.62"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'4.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("D21").Value = "'193.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'140.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -5.01%  "
$ws.Range("D28").Value = "'6.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'15.55"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'0.0491"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'0.898"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "1.112.88"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.803"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'99.21"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").Value = "0.0₆0112"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'55.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +12.29%  "
$ws.Range("D48").Value = "'7.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'0.418"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").Value = "'0.0502"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.41%  "
